# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - Updates the "last updated" timestamp string
#  - Updates case counters for several countries (rows 4, 27, 29, 64, 66, 73, 74)
#  - Georgia overtakes Uruguay in the ranking: row 151 (previously Uruguay)
#    becomes Georgia with its new totals, and row 152 (previously Georgia)
#    becomes Uruguay with the totals Uruguay had before (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp (row 1, column A) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 08:34"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6460421
$ws.Range("C4").Value = 171
$ws.Range("D4").Value = 3726099
$ws.Range("E4").Value = 2541069
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 193253

# --- Row 27: Ucrania ---
$ws.Range("B27").Value = 138068
$ws.Range("C27").Value = 2174
$ws.Range("D27").Value = 62606
$ws.Range("E27").Value = 72585
$ws.Range("G27").Value = 31
$ws.Range("H27").Value = 2877

# --- Row 29: Israel ---
$ws.Range("B29").Value = 131641
$ws.Range("C29").Value = 997
$ws.Range("D29").Value = 103846
$ws.Range("E29").Value = 26776

# --- Row 64: Uzbekistan ---
$ws.Range("B64").Value = 43775
$ws.Range("C64").Value = 188
$ws.Range("E64").Value = 2148
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 350

# --- Row 66: Afganistan ---
$ws.Range("B66").Value = 38494
$ws.Range("C66").Value = 96
$ws.Range("D66").Value = 30557
$ws.Range("E66").Value = 6522
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 1415

# --- Row 73: Australia ---
$ws.Range("B73").Value = 26322
$ws.Range("C73").Value = 43
$ws.Range("D73").Value = 22602
$ws.Range("E73").Value = 2958

# --- Row 74: El Salvador ---
$ws.Range("D74").Value = 15822
$ws.Range("E74").Value = 9722
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 764

# --- Rows 151/152: Georgia overtakes Uruguay ---
# Row 151 keeps the "Georgia" label (shared-string slot 155) but gets the
# country's newly updated figures.
$ws.Range("B151").Value = 1684
$ws.Range("C151").Value = 34
$ws.Range("D151").Value = 1315
$ws.Range("E151").Value = 350
$ws.Range("H151").Value = 19

# Row 152 keeps the "Uruguay" label (shared-string slot 156) with the
# totals Uruguay already had (unchanged from before the update).
$ws.Range("B152").Value = 1679
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 1459
$ws.Range("E152").Value = 175
$ws.Range("H152").Value = 45

# The country names themselves swap order: Georgia now precedes Uruguay.
$ws.Range("A151").Value = "Georgia"
$ws.Range("A152").Value = "Uruguay"
